# "name issue is fixed"
# Bumps the "Envio/Envío Semana 09" paths (and related WorkingStart/EndDate
# constants) forward to "Semana 10", and updates the selected cell on the
# Settings / Constants sheets to match where the user left off editing.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Settings sheet: update the week-09 -> week-10 folder paths
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Settings")

$base        = "/Planeacion/0.Envios TS/2022/10 Envío Semana 10/MÉXICO"
$baseNoAccent = "/Planeacion/0.Envios TS/2022/10 Envio Semana 10/MÉXICO"

$ws1.Range("B2").Value  = "$base/Base de Datos"
$ws1.Range("B6").Value  = "$base/Base de Datos"
$ws1.Range("B8").Value  = "$baseNoAccent/Base de Datos/Consolidado"
$ws1.Range("B10").Value = "$base/VIP"
$ws1.Range("B12").Value = "$base/Tradicional"
$ws1.Range("B14").Value = "$base/Mi Cine"
$ws1.Range("B16").Value = "$base/VIP"
$ws1.Range("B18").Value = "$base/Atmosfera"
$ws1.Range("B20").Value = "$base/Base de Datos"
$ws1.Range("B23").Value = "$base/Base de Datos/Exportadas"
$ws1.Range("B26").Value = "$base/Base de Datos/Exportadas"
$ws1.Range("B29").Value = "$base/Base de Datos/Exportadas"
$ws1.Range("B33").Value = "$base/Base de Datos/Exportadas"

# ---------------------------------------------------------------------
# Constants sheet: bump the working start/end dates forward one week
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Constants")

$ws2.Range("B24").Value = 44621
$ws2.Range("B25").Value = 44626

# ---------------------------------------------------------------------
# Restore the cursor/selection position on each sheet
# ---------------------------------------------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("B6").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("B26").Select() | Out-Null

Write-Output "Applied week 09 -> week 10 config updates"
